# Daily Report update: 2026-01-22
# Adds the new daily snapshot (date serial 46043 = 2026-01-21) as 22 new rows
# at the bottom of Daily_Data, then refreshes the dependent summary rows on
# Today_Summary (MANFRA, TORDELLA & BROOKES, LLC) and Monthly_Stats
# (grand total row + the MANFRA Eligible monthly detail row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Daily_Data: append the new day's rows (266-287)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Daily_Data")

$newDate = 46043

$newRows = @(
    @{B="ASAHI DEPOSITORY LLC Registered"; C=0; D=0; E=0; F=0; G=0; H=0}
    @{B="ASAHI DEPOSITORY LLC Eligible"; C=0; D=0; E=0; F=0; G=0; H=0}
    @{B="BRINK'S, INC. Registered"; C=91733.761; D=0; E=0; F=0; G=0; H=91733.761}
    @{B="BRINK'S, INC. Eligible"; C=27494.288; D=0; E=0; F=0; G=0; H=27494.288}
    @{B="CNT DEPOSITORY, INC. Registered"; C=1246.06; D=0; E=0; F=0; G=0; H=1246.06}
    @{B="CNT DEPOSITORY, INC. Eligible"; C=0; D=0; E=0; F=0; G=0; H=0}
    @{B="DELAWARE DEPOSITORY Registered"; C=1633.941; D=0; E=0; F=0; G=0; H=1633.941}
    @{B="DELAWARE DEPOSITORY Eligible"; C=18459.584; D=0; E=0; F=0; G=0; H=18459.584}
    @{B="HSBC BANK, USA Registered"; C=1394.758; D=0; E=0; F=0; G=0; H=1394.758}
    @{B="HSBC BANK, USA Eligible"; C=9281.978999999999; D=0; E=0; F=0; G=0; H=9281.978999999999}
    @{B="INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"; C=2395.448; D=0; E=0; F=0; G=0; H=2395.448}
    @{B="INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"; C=0; D=0; E=0; F=0; G=0; H=0}
    @{B="JP MORGAN CHASE BANK NA Registered"; C=114985.579; D=0; E=0; F=0; G=0; H=114985.579}
    @{B="JP MORGAN CHASE BANK NA Eligible"; C=135413.823; D=0; E=0; F=0; G=0; H=135413.823}
    @{B="LOOMIS INTERNATIONAL (US) LLC Registered"; C=63745.991; D=0; E=0; F=0; G=0; H=63745.991}
    @{B="LOOMIS INTERNATIONAL (US) LLC Eligible"; C=132077.206; D=0; E=0; F=0; G=0; H=132077.206}
    @{B="MALCA-AMIT USA, LLC Registered"; C=395.145; D=0; E=0; F=0; G=0; H=395.145}
    @{B="MALCA-AMIT USA, LLC Eligible"; C=0; D=0; E=0; F=0; G=0; H=0}
    @{B="MANFRA, TORDELLA & BROOKES, LLC Registered"; C=50220.42; D=0; E=0; F=0; G=0; H=50220.42}
    @{B="MANFRA, TORDELLA & BROOKES, LLC Eligible"; C=11149.237; D=0; E=9877.864; F=-9877.864; G=0; H=1271.373}
    @{B="STONEX PRECIOUS METALS LLC Registered"; C=14122.765; D=0; E=0; F=0; G=0; H=14122.765}
    @{B="STONEX PRECIOUS METALS LLC Eligible"; C=16.075; D=0; E=0; F=0; G=0; H=16.075}
)

$lastRow = $ws.UsedRange.Rows.Count()
$dateStyleSource = $ws.Cells.Item($lastRow, 1)

$r = $lastRow + 1
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 1).NumberFormat = $dateStyleSource.NumberFormat()
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Today_Summary: MANFRA, TORDELLA & BROOKES, LLC eligible/registered totals
#    now reflect the new day's balances.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Today_Summary")
$wsSummary.Cells.Item(11, 2).Value = 1271.373
$wsSummary.Cells.Item(11, 4).Value = 51491.793

# ---------------------------------------------------------------------------
# 3. Monthly_Stats: month-to-date grand total and the MANFRA Eligible detail
#    row both shift by the -9877.864 withdrawal recorded today.
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")
$wsMonthly.Cells.Item(2, 2).Value = 324014.328
$wsMonthly.Cells.Item(2, 4).Value = 665888.196

$wsMonthly.Cells.Item(25, 4).Value = 37704.598
$wsMonthly.Cells.Item(25, 5).Value = 1271.373
